# Added Test scripts PUBLONS028.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Populate the "Results" column (D) with "Y" for rows 28, 29, 30, 32, 33
$ws.Range("D28").Value = "Y"
$ws.Range("D29").Value = "Y"
$ws.Range("D30").Value = "Y"
$ws.Range("D32").Value = "Y"
$ws.Range("D33").Value = "Y"

# Update the view state: scrolled/selected position
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("C37").Select()
